$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.179.58'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.854.97'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.86'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6890'
$ws.Range('E6').Value = '  -1.27%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07778'
$ws.Range('E8').Value = '  +4.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3047'
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('E10').Value = '  -2.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08067'
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.862.20'
$ws.Range('E12').Value = '  -2.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7213'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.188'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.27'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.179.17'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.737'
$ws.Range('E17').Value = '  -2.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007808'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.25'
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '234.90'
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9995'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('E22').Value = '  -1.24%  '
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.461'
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.07'
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.966'
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('E27').Value = '  -4.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.03'
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.399'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.513'
$ws.Range('E31').Value = '  +2.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.487'
$ws.Range('E32').Value = '  -1.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.007'
$ws.Range('E33').Value = '  -1.43%  '
$ws.Range('E34').Value = '  -1.47%  '
$ws.Range('E35').Value = '  -1.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7024'
$ws.Range('E36').Value = '  -2.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9997'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01846'
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9390'
$ws.Range('E41').Value = '  +7.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.093.22'
$ws.Range('E42').Value = '  +5.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.958'
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4286'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.44'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9998'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.33'
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.798'
$ws.Range('E48').Value = '  +2.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.007.14'
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.159'
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.006'
$ws.Range('E51').Value = '  -3.57%  '
